$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.755.80'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.637.84'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'607.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').Value = "'147.14"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.09%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = "'0.589"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').Value = "'0.383"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.37%  '
$ws.Range('D11').Value = "'5.59"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').Value = "'27.37"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').Value = '3.112.40'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = '63.605.02'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '2.632.47'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('D18').Value = "'11.74"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('D20').Value = "'346.54"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Value = "'5.56"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').Value = "'66.30"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('E25').Value = '  +7.13%  '
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').Value = "'9.24"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.17%  '
$ws.Range('D28').Value = "'563.45"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.46%  '
$ws.Range('D29').Value = "'8.10"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.54%  '
$ws.Range('D30').Value = "'0.161"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('D31').Value = "'0.999"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').Value = '0.0₃0852'
$ws.Range('E33').Value = '  +4.78%  '
$ws.Range('E34').Value = '  -2.15%  '
$ws.Range('D35').Value = "'5.28"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.85%  '
$ws.Range('D36').Value = "'169.47"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = "'0.404"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  +5.03%  '
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D42').Value = "'164.93"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.78%  '
$ws.Range('D43').Value = "'40.15"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = "'3.79"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('D45').Value = "'21.86"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.13%  '
$ws.Range('D46').Value = "'0.0565"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = "'0.625"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('E48').Value = '  +14.00%  '
$ws.Range('D49').Value = "'0.0244"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.33%  '
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('E51').Value = '  -1.49%  '
